# Added marks of Made Easy GATE Mock 1
# The author cleared a batch of "weak topic" cells (now mastered / tracked
# elsewhere) from the tracker on Sheet1, and left the selection on D13.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear the topic cells that were removed from the tracker.
$cellsToClear = @(
    "C4", "D4", "E4", "F4", "G4", "H4", "I4", "J4", "L4", "M4", "N4",
    "C10", "H10",
    "G11", "H11", "K11", "L11", "M11"
)

foreach ($ref in $cellsToClear) {
    $ws.Range($ref).ClearContents()
}

# Restore normal view (no frozen/scrolled topLeftCell) and move the active
# selection to D13, matching the saved workbook state.
$ws.Activate()
$ws.Range("D13").Select()
